$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186790347099304
$ws.Range("B1").Value = 1.17013418674469
$ws.Range("C1").Value = 6.788376808166504
$ws.Range("D1").Value = 2.065802574157715
$ws.Range("E1").Value = 1.139171838760376
